# Atualizacao de bases das ligas, do dia: 28-05-2024 as 07:50
#
# The source rows had been written with swapped/rotated data relative to
# their match id (column B) and the rest of the row (C..AD). This script
# restores the correct pairing by moving each row's B..AD content to the
# row it actually belongs to:
#   - row 93  <-> row 94   (full swap)
#   - row 95  -> row 96 -> row 99 -> row 95   (3-way rotation)
#   - row 114 <-> row 115  (full swap)
# Column A (the running index) is left untouched throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # B
$lastCol  = 30  # AD

function Get-RowValues($ws, $row, $firstCol, $lastCol) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += , ($ws.Cells.Item($row, $c).Value())
    }
    return $vals
}

function Set-RowValues($ws, $row, $firstCol, $lastCol, $vals) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$i]
        $i = $i + 1
    }
}

# Snapshot every row involved before any writes happen, so later writes
# never read already-overwritten data.
$row93 = Get-RowValues $ws 93 $firstCol $lastCol
$row94 = Get-RowValues $ws 94 $firstCol $lastCol
$row95 = Get-RowValues $ws 95 $firstCol $lastCol
$row96 = Get-RowValues $ws 96 $firstCol $lastCol
$row99 = Get-RowValues $ws 99 $firstCol $lastCol
$row114 = Get-RowValues $ws 114 $firstCol $lastCol
$row115 = Get-RowValues $ws 115 $firstCol $lastCol

# 93 <-> 94
Set-RowValues $ws 93 $firstCol $lastCol $row94
Set-RowValues $ws 94 $firstCol $lastCol $row93

# 95 -> 96 -> 99 -> 95
Set-RowValues $ws 95 $firstCol $lastCol $row96
Set-RowValues $ws 96 $firstCol $lastCol $row99
Set-RowValues $ws 99 $firstCol $lastCol $row95

# 114 <-> 115
Set-RowValues $ws 114 $firstCol $lastCol $row115
Set-RowValues $ws 115 $firstCol $lastCol $row114
